$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The schedule shifts: a new "Android Compose betanulás" task is inserted
# for weeks 7-9 (rows 8-10), pushing the existing Android client screens
# work (kezdőképernyő / termékek-kosár-checkout+filterezés) down into
# weeks 10-11 (rows 11-12). Thesis writing ("Diplomamunka írás") still
# occupies weeks 12-14 (rows 13-15), unchanged.

$ws.Range("C8").Value = "Android Compose betanulás"
$ws.Range("C9").Value = "Android Compose betanulás"
$ws.Range("C10").Value = "Android Compose betanulás"
$ws.Range("C11").Value = "Android kliens - kezdőképernyő, autentikáció megvalósítás, lokális adatbázis megvalósítás"
$ws.Range("C12").Value = "Android kliens - termékek, kosár, checkout képernyők megvalósítása,  termék filterezés megvalósítása"

# Move the active selection to C14 (matches the author's cursor position
# at save time).
$ws.Range("C14").Select()
